# Insert a new weekly price record at row 31 for
# "Vega Monumental Concepción - Jengibre", pushing the existing
# rows 31-38 down to 32-39 (dimension grows from A1:R38 to A1:R39).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 31, shifting rows 31:38 down to 32:39.
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new record.
$ws.Range("A31").Value = 11
$ws.Range("B31").Value = "Vega Monumental Concepción"
$ws.Range("C31").Value = "Bíobío"
$ws.Range("D31").Value = 44782
$ws.Range("E31").Value = 8
$ws.Range("F31").Value = 100114007
$ws.Range("G31").Value = "Jengibre"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 40
$ws.Range("K31").Value = 13000
$ws.Range("L31").Value = 14000
$ws.Range("M31").Value = 13500
$ws.Range("N31").Value = "`$/caja 13 kilos"
$ws.Range("O31").Value = "Perú"
$ws.Range("P31").Value = 1038
$ws.Range("Q31").Value = 13
$ws.Range("R31").Value = "Hortaliza"
